# Improvements made to English localization of text translated from Italian
# Update gendered pronouns ("him"/"his"/"he") to gender-neutral language
# ("them"/"their"/"they"/"others") in the en-gb column of the
# QuestionLocalizations sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuestionLocalizations")

$ws.Range("B5").Value = "I attempt to deal with all of their and my concerns."
$ws.Range("B25").Value = "I will let others have some of their positions if they lets me have some of mine."
$ws.Range("B28").Value = "I tell them my ideas and ask them for theirs."
$ws.Range("B29").Value = "I try to show them the logic and benefits of my position."
$ws.Range("B36").Value = "If it makes the other person happy, I might let them maintain their views."
$ws.Range("B37").Value = "I will let them have some of their positions if they lets me have some of mine."
$ws.Range("B44").Value = "I try to find a position that is intermediate between theirs and mine."
$ws.Range("B48").Value = "If the other's position seems very important to them, I would try to meet their wishes."
$ws.Range("B49").Value = "I try to get others to settle for a compromise solution."
$ws.Range("B50").Value = "I try to show others the logic and benefits of my position."
$ws.Range("B55").Value = "If it makes the other person happy, I might let them maintain their views."
